# Update the betting-odds figures for the 4 matches that changed in the
# 2025-06-06 FlashScore refresh (rows 2, 3, 8 and 14 of Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Def. de Belgrano - Almirante Brown)
$ws.Range("J2").Value = 1.13
$ws.Range("K2").Value = 6
$ws.Range("N2").Value = 3.1
$ws.Range("O2").Value = 1.36

# Row 3 (Los Andes - Alvarado)
$ws.Range("G3").Value = 1.8
$ws.Range("I3").Value = 5.75
$ws.Range("J3").Value = 1.17
$ws.Range("K3").Value = 5
$ws.Range("N3").Value = 3.4
$ws.Range("O3").Value = 1.33
$ws.Range("R3").Value = 2.75
$ws.Range("S3").Value = 1.4
$ws.Range("AD3").Value = 9.5

# Row 8 (SJK Akatemia - JIPPO)
$ws.Range("G8").Value = 3.2
$ws.Range("I8").Value = 2.1
$ws.Range("M8").Value = 3.4
$ws.Range("O8").Value = 1.93
$ws.Range("T8").Value = 11.75
$ws.Range("U8").Value = 19
$ws.Range("V8").Value = 11
$ws.Range("W8").Value = 45
$ws.Range("X8").Value = 26
$ws.Range("Y8").Value = 28
$ws.Range("Z8").Value = 11.5
$ws.Range("AA8").Value = 6.6
$ws.Range("AE8").Value = 11.25
$ws.Range("AG8").Value = 21
$ws.Range("AH8").Value = 16

# Row 14 (Hartford Athletic - North Carolina)
$ws.Range("G14").Value = 3.05
$ws.Range("H14").Value = 3.25
$ws.Range("I14").Value = 2.15
$ws.Range("J14").Value = 1.06
$ws.Range("K14").Value = 7.4
$ws.Range("L14").Value = 1.28
$ws.Range("M14").Value = 3.35
$ws.Range("N14").Value = 1.83
$ws.Range("O14").Value = 1.87
$ws.Range("P14").Value = 1.4
$ws.Range("Q14").Value = 2.7
$ws.Range("T14").Value = 10.25
$ws.Range("U14").Value = 17
$ws.Range("V14").Value = 10.75
$ws.Range("W14").Value = 40
$ws.Range("X14").Value = 26
$ws.Range("Y14").Value = 32
$ws.Range("Z14").Value = 7.4
$ws.Range("AA14").Value = 6.4
$ws.Range("AD14").Value = 8.25
$ws.Range("AE14").Value = 11
$ws.Range("AF14").Value = 8.75
$ws.Range("AG14").Value = 22
$ws.Range("AH14").Value = 17
$ws.Range("AI14").Value = 25
